$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text corrections
$ws.Range("B2").Value  = "<kero>"
$ws.Range("B5").Value  = "<long>"
$ws.Range("B7").Value  = "<upa>"
$ws.Range("B9").Value  = "<nomber>"
$ws.Range("B10").Value = "<four>"
$ws.Range("B11").Value = "<would>"
$ws.Range("B13").Value = "<bac>"
$ws.Range("B15").Value = "<calter>"

# Column C numeric corrections
$ws.Range("C2").Value  = 22
$ws.Range("C3").Value  = 26
$ws.Range("C4").Value  = 27
$ws.Range("C6").Value  = 28
$ws.Range("C7").Value  = 31
$ws.Range("C8").Value  = 33
$ws.Range("C9").Value  = 21
$ws.Range("C10").Value = 30
$ws.Range("C12").Value = 33
$ws.Range("C14").Value = 31
$ws.Range("C15").Value = 8
